$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 7 data (new parish: Bussigny - Villars-Sainte-Croix)
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = 3050000000
$ws.Range("C7").Value = 3050
$ws.Range("D7").Value = 3000
$ws.Range("E7").Value = "Bussigny – Villars-Sainte-Croix"
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = "P"

# Match the styling of the preceding data rows (style index 4: vertical-top alignment)
# Only the numeric columns (A-D, F) carry the style; E/G (text) stay default, like rows 2-6.
$ws.Range("A7:D7").VerticalAlignment = $ws.Range("A6:D6").VerticalAlignment
$ws.Range("F7").VerticalAlignment = $ws.Range("F6").VerticalAlignment

# Update the selection to mirror the post-edit state
$ws.Range("E9").Select()
